# Applies the scheduled Universalis market-data refresh to the per-job leve
# profit tables (columns H:N = currentAveragePrice*, LevePrice*, LeveProfit*)
# across all eight Marilith_Profits job sheets.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 167.66667
$ws.Range("I33").Value = 154.64285
$ws.Range("J33").Value = 350
$ws.Range("K33").Value = 154.64285
$ws.Range("L33").Value = 350
$ws.Range("M33").Value = 74.35714999999999
$ws.Range("N33").Value = -808
$ws.Range("H69").Value = 6916
$ws.Range("I69").Value = 1000
$ws.Range("J69").Value = 8099.2
$ws.Range("K69").Value = 3000
$ws.Range("L69").Value = 24297.6
$ws.Range("M69").Value = -2126
$ws.Range("N69").Value = -26045.6
$ws.Range("H72").Value = 6916
$ws.Range("I72").Value = 1000
$ws.Range("J72").Value = 8099.2
$ws.Range("K72").Value = 9000
$ws.Range("L72").Value = 72892.8
$ws.Range("M72").Value = -4632
$ws.Range("N72").Value = -81628.8
$ws.Range("H74").Value = 3495.5
$ws.Range("I74").Value = 3495.5
$ws.Range("K74").Value = 3495.5
$ws.Range("M74").Value = -2559.5
$ws.Range("H77").Value = 3495.5
$ws.Range("I77").Value = 3495.5
$ws.Range("K77").Value = 17477.5
$ws.Range("M77").Value = -12797.5
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents() | Out-Null
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents() | Out-Null
$ws.Range("H92").Value = 384.2
$ws.Range("I92").Value = 343.30768
$ws.Range("K92").Value = 343.30768
$ws.Range("M92").Value = 904.69232
$ws.Range("H99").Value = 752.5454999999999
$ws.Range("I99").Value = 612.1429000000001
$ws.Range("K99").Value = 1836.4287
$ws.Range("M99").Value = -338.4287000000002

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H55").Value = 18997.223
$ws.Range("J55").Value = 18997.223
$ws.Range("L55").Value = 18997.223
$ws.Range("N55").Value = -19627.223
$ws.Range("H80").Value = 29999.166
$ws.Range("J80").Value = 29999.166
$ws.Range("L80").Value = 29999.166
$ws.Range("N80").Value = -31995.166
$ws.Range("H83").Value = 29999.166
$ws.Range("J83").Value = 29999.166
$ws.Range("L83").Value = 89997.49800000001
$ws.Range("N83").Value = -99981.49800000001

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 15000
$ws.Range("J35").Value = 15000
$ws.Range("L35").Value = 15000
$ws.Range("N35").Value = -15620
$ws.Range("H82").Value = 19539.525
$ws.Range("I82").Value = 12296
$ws.Range("J82").Value = 29499.375
$ws.Range("K82").Value = 12296
$ws.Range("L82").Value = 29499.375
$ws.Range("M82").Value = -11913
$ws.Range("N82").Value = -30265.375
$ws.Range("H85").Value = 19539.525
$ws.Range("I85").Value = 12296
$ws.Range("J85").Value = 29499.375
$ws.Range("K85").Value = 12296
$ws.Range("L85").Value = 29499.375
$ws.Range("M85").Value = -10970
$ws.Range("N85").Value = -32151.375
$ws.Range("H105").Value = 2749.8
$ws.Range("I105").Value = 2642.5715
$ws.Range("J105").Value = 3000
$ws.Range("K105").Value = 2642.5715
$ws.Range("L105").Value = 3000
$ws.Range("M105").Value = -895.5715
$ws.Range("N105").Value = -6494

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 2366
$ws.Range("J15").Value = 2049
$ws.Range("L15").Value = 2049
$ws.Range("N15").Value = -2389
$ws.Range("H41").Value = 14124.75
$ws.Range("J41").Value = 15714
$ws.Range("L41").Value = 15714
$ws.Range("N41").Value = -16570
$ws.Range("H50").Value = 19462.924
$ws.Range("J50").Value = 18665.834
$ws.Range("L50").Value = 18665.834
$ws.Range("N50").Value = -19915.834
$ws.Range("H60").Value = 22298.715
$ws.Range("H96").Value = 23262
$ws.Range("J96").Value = 23262
$ws.Range("L96").Value = 23262
$ws.Range("N96").Value = -28754
$ws.Range("H132").Value = 3718.16
$ws.Range("I132").Value = 2512.5715
$ws.Range("J132").Value = 4187
$ws.Range("K132").Value = 7537.7145
$ws.Range("L132").Value = 12561
$ws.Range("M132").Value = -5007.7145
$ws.Range("N132").Value = -17621

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 1714.8334
$ws.Range("I25").Value = 924.75
$ws.Range("J25").Value = 3295
$ws.Range("K25").Value = 2774.25
$ws.Range("L25").Value = 9885
$ws.Range("M25").Value = -2605.25
$ws.Range("N25").Value = -10223
$ws.Range("H30").Value = 1714.8334
$ws.Range("I30").Value = 924.75
$ws.Range("J30").Value = 3295
$ws.Range("K30").Value = 2774.25
$ws.Range("L30").Value = 9885
$ws.Range("M30").Value = -2672.25
$ws.Range("N30").Value = -10089
$ws.Range("H39").Value = 2148.5
$ws.Range("J39").Value = 2148.5
$ws.Range("L39").Value = 6445.5
$ws.Range("N39").Value = -7033.5
$ws.Range("H46").Value = 444
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents() | Out-Null
$ws.Range("H122").Value = 574.25
$ws.Range("I122").Value = 548.2
$ws.Range("J122").Value = 617.6667
$ws.Range("K122").Value = 4933.8
$ws.Range("L122").Value = 5559.0003
$ws.Range("M122").Value = -2483.8
$ws.Range("N122").Value = -10459.0003
$ws.Range("H125").Value = 3500
$ws.Range("I125").Value = 3500
$ws.Range("K125").Value = 10500
$ws.Range("M125").Value = -5580

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 11500
$ws.Range("J54").Value = 11500
$ws.Range("L54").Value = 11500
$ws.Range("N54").Value = -12280
$ws.Range("H57").Value = 19998
$ws.Range("J57").Value = 19998
$ws.Range("L57").Value = 19998
$ws.Range("N57").Value = -21638

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 8800.799999999999
$ws.Range("J19").Value = 8800.799999999999
$ws.Range("L19").Value = 8800.799999999999
$ws.Range("N19").Value = -9140.799999999999
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents() | Out-Null
$ws.Range("H136").Value = 3701.1667
$ws.Range("I136").Value = 3114.25
$ws.Range("K136").Value = 9342.75
$ws.Range("M136").Value = -6792.75

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1013333.3
$ws.Range("J54").Value = 1013333.3
$ws.Range("L54").Value = 1013333.3
$ws.Range("N54").Value = -1014373.3
$ws.Range("H122").Value = 1900
$ws.Range("I122").Value = 300
$ws.Range("K122").Value = 900
$ws.Range("M122").Value = 1550
$ws.Range("H136").Value = 5217.6
$ws.Range("I136").Value = 1856.5714
$ws.Range("K136").Value = 5569.7142
$ws.Range("M136").Value = -3019.7142
